# Actualización automática 2025-11-06 17:30:07
#
# A new client "ECUAFERRI S.A." is inserted alphabetically (between
# "DANIELA ELIZABETH BECERRA BECERRA" and "EQUISAB S.A.") as row 21 in
# both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing every
# following row down by one. The trailing summary row in each sheet
# moves down accordingly, and on "VENTAS POR GRUPO" the "X de 51" counts
# in that summary row become "X de 52" to reflect the extra client row.

$wb = $excel.ActiveWorkbook

# ---- Sheet: VENTAS POR GRUPO (18 columns, A:R) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(21).Insert()
$ws1.Cells.Item(21, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(21, 2).Value = "ECUAFERRI S.A."
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(21, $c).Value = 0
}

# The old last row (53 -> now 54) holds "X de 51" labels; bump the count
# to reflect the newly added client row.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(54, $c)
    $text = $cell.Value2
    $cell.Value = $text -replace "de 51", "de 52"
}

# ---- Sheet: VENTA MENSUAL (7 columns, A:G) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(21).Insert()
$ws2.Cells.Item(21, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(21, 2).Value = "ECUAFERRI S.A."
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(21, $c).Value = 0
}
